$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new data rows right before the current row 435, pushing the
# existing rows 435-470 down to 438-473 (new dimension A1:R473).
$ws.Rows.Item(435).EntireRow.Insert()
$ws.Rows.Item(435).EntireRow.Insert()
$ws.Rows.Item(435).EntireRow.Insert()

# Populate the new row 435 (new weekly entry, "Primera" quality, $/caja 36 atados).
$ws.Range("A435").Value = 9
$ws.Range("B435").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C435").Value = "Metropolitana"
$ws.Range("D435").Value = 44578
$ws.Range("E435").Value = 13
$ws.Range("F435").Value = 100112040
$ws.Range("G435").Value = "Cilantro"
$ws.Range("H435").Value = "Sin especificar"
$ws.Range("I435").Value = "Primera"
$ws.Range("J435").Value = 43
$ws.Range("K435").Value = 8000
$ws.Range("L435").Value = 8000
$ws.Range("M435").Value = 8000
$ws.Range("N435").Value = "$/caja 36 atados"
$ws.Range("O435").Value = "Región Metropolitana"
$ws.Range("P435").Value = 222
$ws.Range("Q435").Value = 36
$ws.Range("R435").Value = "Hortaliza"

# Populate the new row 436 (new weekly entry, "Primera" quality, $/docena de atados).
$ws.Range("A436").Value = 9
$ws.Range("B436").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C436").Value = "Metropolitana"
$ws.Range("D436").Value = 44578
$ws.Range("E436").Value = 13
$ws.Range("F436").Value = 100112040
$ws.Range("G436").Value = "Cilantro"
$ws.Range("H436").Value = "Sin especificar"
$ws.Range("I436").Value = "Primera"
$ws.Range("J436").Value = 97
$ws.Range("K436").Value = 16000
$ws.Range("L436").Value = 18000
$ws.Range("M436").Value = 16990
$ws.Range("N436").Value = "$/docena de atados"
$ws.Range("O436").Value = "Región Metropolitana"
$ws.Range("P436").Value = 5663
$ws.Range("Q436").Value = 3
$ws.Range("R436").Value = "Hortaliza"

# Populate the new row 437 (new weekly entry, "Segunda" quality, $/docena de atados).
$ws.Range("A437").Value = 9
$ws.Range("B437").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C437").Value = "Metropolitana"
$ws.Range("D437").Value = 44578
$ws.Range("E437").Value = 13
$ws.Range("F437").Value = 100112040
$ws.Range("G437").Value = "Cilantro"
$ws.Range("H437").Value = "Sin especificar"
$ws.Range("I437").Value = "Segunda"
$ws.Range("J437").Value = 43
$ws.Range("K437").Value = 12000
$ws.Range("L437").Value = 13000
$ws.Range("M437").Value = 12488
$ws.Range("N437").Value = "$/docena de atados"
$ws.Range("O437").Value = "Región Metropolitana"
$ws.Range("P437").Value = 4163
$ws.Range("Q437").Value = 3
$ws.Range("R437").Value = "Hortaliza"
